# The "area" / QUALITY / 7 entry (row 37) was removed from the gold-label
# sheet. Deleting the entire row shifts every subsequent row up by one,
# which matches the rest of the diff (row 38 -> 37, row 39 -> 38, ...,
# row 271 -> 270) without needing to touch any other cell individually.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("A37").EntireRow.Delete()
